$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the "Good Drivers" table: a newly-observed driver
# version (21.40.1.3) is added at the top of the list, pushing the
# existing rows (and the trailing blank padding rows) down by one.
$ws.Rows.Item(12).Insert()

# New row 12: the newly observed driver. It has no "Driver Vintage" yet.
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").Value = 11128
$ws.Range("B12").NumberFormat = "#,##0"
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 0

# This week's updated sample counts for the existing drivers (the insert
# above shifted the old rows 12-17 down to 13-18; only the counts moved).
$ws.Range("B13").Value = 486214
$ws.Range("B14").Value = 79953
$ws.Range("B15").Value = 35355
$ws.Range("B16").Value = 65425
$ws.Range("B17").Value = 117653
